$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = strikeouts) values per regenerated save_data
$ws.Range("G2").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 2
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 2
$ws.Range("G20").Value = 1
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("G24").Value = 2
$ws.Range("G25").Value = 2
$ws.Range("G26").Value = 0
$ws.Range("G27").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("G30").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("G32").Value = 0
$ws.Range("G33").Value = 2
$ws.Range("G34").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("G36").Value = 2
$ws.Range("G37").Value = 3
$ws.Range("G38").Value = 0
$ws.Range("G39").Value = 2
$ws.Range("G40").Value = 2
$ws.Range("G41").Value = 1
$ws.Range("G42").Value = 2
$ws.Range("G43").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("G45").Value = 1
$ws.Range("G46").Value = 1
$ws.Range("G47").Value = 4
$ws.Range("G48").Value = 1
$ws.Range("G49").Value = 1
$ws.Range("G50").Value = 0
$ws.Range("G51").Value = 2
$ws.Range("G52").Value = 1
$ws.Range("G53").Value = 2
$ws.Range("G54").Value = 2
$ws.Range("G55").Value = 1
$ws.Range("G56").Value = 1
$ws.Range("G58").Value = 2
$ws.Range("G59").Value = 1
$ws.Range("G60").Value = 2
$ws.Range("G61").Value = 2
$ws.Range("G62").Value = 1
$ws.Range("G63").Value = 1
$ws.Range("G66").Value = 3

Write-Host "Updated G column values."
